$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing data down
$ws.Rows.Item(1).Insert()

# Add header labels
$ws.Range("A1").Value = "Ticker"
$ws.Range("B1").Value = "weight"

# Update selection to match the target state
$ws.Range("D5").Select()
